$d = $word.ActiveDocument

# Locate the two target paragraphs by their current text content so the
# script does not depend on a fixed paragraph index.
$groupNumberPara = $null
$groupMembersPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Group Number: 21479*" -and $groupNumberPara -eq $null) {
        $groupNumberPara = $p
    }
    if ($t -like "Group Members: Connor Hay, Alexis Effenberger, Ben Heick, Camryn Rogers, Julianne Heine*" -and $groupMembersPara -eq $null) {
        $groupMembersPara = $p
    }
}

# --- Edit 1: paragraph "Group Number: 21479" -> append a separate run "_3" ---
if ($groupNumberPara -ne $null) {
    $r2 = $groupNumberPara.Range.Duplicate
    $xml2 = '<w:p w:rsidR="00587799" w:rsidRDefault="00587799" w:rsidP="00587799">' + `
              '<w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:jc w:val="center"/></w:pPr>' + `
              '<w:r><w:t>Group Number: 21479</w:t></w:r>' + `
              '<w:r><w:t>_3</w:t></w:r>' + `
            '</w:p>'
    $r2.InsertXML($xml2)
}

# --- Edit 2: paragraph "Group Members: ..." -> split into runs with proofErr markers ---
if ($groupMembersPara -ne $null) {
    $r3 = $groupMembersPara.Range.Duplicate
    $xml3 = '<w:p w:rsidR="00587799" w:rsidRDefault="00587799" w:rsidP="00587799">' + `
              '<w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:jc w:val="center"/></w:pPr>' + `
              '<w:r><w:t xml:space="preserve">Group Members: Connor Hay, Alexis </w:t></w:r>' + `
              '<w:proofErr w:type="spellStart"/>' + `
              '<w:r><w:t>Effenberger</w:t></w:r>' + `
              '<w:proofErr w:type="spellEnd"/>' + `
              '<w:r><w:t xml:space="preserve">, Ben </w:t></w:r>' + `
              '<w:proofErr w:type="spellStart"/>' + `
              '<w:r><w:t>Heick</w:t></w:r>' + `
              '<w:proofErr w:type="spellEnd"/>' + `
              '<w:r><w:t>, Camryn Rogers, Julianne Heine</w:t></w:r>' + `
            '</w:p>'
    $r3.InsertXML($xml3)
}
